# Adapt column header formatting to respective input file names (#7):
#   <suffix>_old -> <suffix>_FV2404
#   <suffix>_new -> <suffix>_FV2410
# and wrap the sheet's used range in an Excel Table, plus freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

# Columns A-J: "<name>_old" -> "<name>_FV2404"
for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $name = $oldHeaders[$i] -replace "_old$", "_FV2404"
    $ws.Cells.Item(1, $i + 1).Value = $name
}

# Column K ("diff") is untouched.

# Columns L-U: "<name>_new" -> "<name>_FV2410"
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $name = $newHeaders[$i] -replace "_new$", "_FV2410"
    $ws.Cells.Item(1, 12 + $i).Value = $name
}

# Turn the whole used range into a native Excel Table ("Table1").
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U84"), 0, 1)
$tbl.Name = "Table1"

# Freeze the header row (split above row 2, freeze).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
